$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the last existing year columns (P) onto the new column Q,
# then fill in the 2020 figures.
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 70.3

# Scroll the view so column C is the left-most visible column, then select Q8
# (matching the saved view state).
$excel.ActiveWindow.SplitColumn = 2
$ws.Range("Q8").Select()
$excel.ActiveWindow.SplitColumn = 0
